# Generate Report for Handoff
# Updates the localization-status workbook: the status cells that read
# "Handed back: in sync with en-US" become "Ready for handoff", and the
# associated timestamps are refreshed to reflect the new handoff-report
# generation time.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns + generate-date column
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-28 22:58:50"

# zh-cn sheet: Status + Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-28 22:58:46"

# de-de sheet: Status (shares the refreshed generate-date text with Overview G2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-28 22:58:50"

# Column widths on the Status columns shrink now that the shorter
# "Ready for handoff" text no longer needs the extra room (closest value
# this host's ColumnWidth quantization can reach to the recorded
# 17.2159881591797 target).
$newStatusColWidth = 98 / 6
$overview.Columns("E:F").ColumnWidth = $newStatusColWidth
$zhcn.Columns("C:C").ColumnWidth = $newStatusColWidth
$dede.Columns("C:C").ColumnWidth = $newStatusColWidth
